# Quarterly rolling-window update: drop the oldest period (column D) and
# append the newest period (new column M), matching the published
# dollar_cumulative income-statement refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the oldest period column (D) and shift everything left.
#    This automatically re-letters all the remaining period columns
#    (old E -> D, old F -> E, ... old M -> L) and drops the now-unused
#    shared strings for the period header / report-date labels that
#    only the old column D referenced.
$ws.Range("D1:D28").Delete()

# 2) Bring the new (now-blank) column M up to the same formatting as the
#    column to its left (L) before filling in values.
$ws.Range("L8:L27").Copy()
$ws.Range("M8:M27").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) The report-date row: the period that used to be reported as
#    "1401-10-29 (6)" (old column J, now shifted to I) has since been
#    reissued as revision 7, dated 1402-02-27.
$ws.Range("I9").Value = "1402-02-27 (7)"

# 4) New column headers for the freshly added period.
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-27"

# 5) New column data for the freshly added period.
$ws.Range("M11").Value = 1288836
$ws.Range("M12").Value = -1178492
$ws.Range("M13").Value = 110344
$ws.Range("M14").Value = -9169
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 16235
$ws.Range("M17").Value = 117410
$ws.Range("M18").Value = -3663
$ws.Range("M19").Value = 17871
$ws.Range("M20").Value = 131617
$ws.Range("M21").Value = -12695
$ws.Range("M22").Value = 118922
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 118922
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 22670
$ws.Range("M27").Value = 0
